$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row with Student ID 545456 (logged at 07:41:13) was removed from the
# scanner log. Deleting the entire row shifts the remaining two log
# entries (122434 / 07:41:18 and 121234 / 07:41:20) up into rows 2-3,
# which also shrinks the sheet's used range from A1:F4 to A1:F3.
$ws.Rows.Item(2).Delete()

# The sheet itself was renamed from "Scanner" to "Session".
$ws.Name = "Session"
